# Commit: "putting in this one data I did before my laptop kills me"
# The author added one more measurement row (row 39) to the File Transfer
# Test Results sheet: a "PC_Local_to_PI_Peter_Jamieson_half_png" run, with
# its Power-While-Running / Power-While-Not-Running readings and the
# existing "= B - C" power-bump formula carried down one more row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy row 38's formatting down to the new row 39 first so the label cell
# picks up the same bold style (s="1") as every other entry in column A.
$ws.Range("A38:D38").Copy($ws.Range("A39:D39"))

# New data row.
$ws.Range("A39").Value = "PC_Local_to_PI_Peter_Jamieson_half_png"
$ws.Range("B39").Value = 3.5703999999999998
$ws.Range("C39").Value = 3.3765000000000001
$ws.Range("D39").Formula = "=B39-C39"

# Leave the cursor where the author left it: on the freshly typed D39.
$ws.Range("D39").Select() | Out-Null
